$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.010.62"
$ws.Range("E2").Value = "  +2.21%  "
$ws.Range("D3").Value = "'1.675.90"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("D5").Value = "'329.77"
$ws.Range("E5").Value = "  +7.33%  "
$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").Value = "'0.3659"
$ws.Range("E7").Value = "  +1.35%  "
$ws.Range("D8").Value = "'47.14"
$ws.Range("E8").Value = "  -0.61%  "
$ws.Range("D9").Value = "'0.3240"
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("D10").Value = "'1.149"
$ws.Range("E10").Value = "  +2.81%  "
$ws.Range("D11").Value = "'0.07186"
$ws.Range("E11").Value = "  +2.90%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").Value = "'6.095"
$ws.Range("E13").Value = "  +3.66%  "
$ws.Range("D14").Value = "'19.70"
$ws.Range("E14").Value = "  +1.69%  "
$ws.Range("D15").Value = "'1.670.34"
$ws.Range("E15").Value = "  +0.83%  "
$ws.Range("D16").Value = "'6.669"
$ws.Range("E16").Value = "  +1.67%  "
$ws.Range("E17").Value = "  +0.75%  "
$ws.Range("D18").Value = "'0.06539"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").Value = "'0.9993"
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").Value = "'78.97"
$ws.Range("E20").Value = "  +3.52%  "
$ws.Range("D21").Value = "'15.85"
$ws.Range("E21").Value = "  +1.45%  "
$ws.Range("D22").Value = "'5.920"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").Value = "'12.97"
$ws.Range("E23").Value = "  +3.54%  "
$ws.Range("D24").Value = "'24.991.43"
$ws.Range("E24").Value = "  +2.13%  "
$ws.Range("E25").Value = "  -1.25%  "
$ws.Range("D26").Value = "'2.380"
$ws.Range("E26").Value = "  +2.89%  "
$ws.Range("D27").Value = "'149.35"
$ws.Range("E27").Value = "  +1.75%  "
$ws.Range("D28").Value = "'18.73"
$ws.Range("E28").Value = "  +1.67%  "
$ws.Range("D29").Value = "'1.857.75"
$ws.Range("E29").Value = "  +0.95%  "
$ws.Range("D30").Value = "'126.13"
$ws.Range("E30").Value = "  +1.95%  "
$ws.Range("D31").Value = "'1.200"
$ws.Range("E31").Value = "  +2.66%  "
$ws.Range("D32").Value = "'4.091"
$ws.Range("E32").Value = "  +3.02%  "
$ws.Range("D33").Value = "'5.810"
$ws.Range("E33").Value = "  +3.56%  "
$ws.Range("D34").Value = "'0.08464"
$ws.Range("E34").Value = "  +0.96%  "
$ws.Range("D35").Value = "'1.671"
$ws.Range("E35").Value = "  -1.61%  "
$ws.Range("D36").Value = "'12.35"
$ws.Range("E36").Value = "  +0.41%  "
$ws.Range("D37").Value = "'5.166"
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("D38").Value = "'0.06094"
$ws.Range("E38").Value = "  +0.89%  "
$ws.Range("E39").Value = "  +2.56%  "
$ws.Range("D40").Value = "'0.2093"
$ws.Range("E40").Value = "  +2.22%  "
$ws.Range("D41").Value = "'0.02231"
$ws.Range("E41").Value = "  +1.32%  "
$ws.Range("D42").Value = "'8.270"
$ws.Range("E42").Value = "  +0.78%  "
$ws.Range("D43").Value = "'0.9993"
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("D44").Value = "'0.5963"
$ws.Range("E44").Value = "  +1.28%  "
$ws.Range("D45").Value = "'13.66"
$ws.Range("E45").Value = "  +8.64%  "
$ws.Range("D46").Value = "'3.831"
$ws.Range("E46").Value = "  +2.54%  "
$ws.Range("D47").Value = "'0.5731"
$ws.Range("E47").Value = "  +2.80%  "
$ws.Range("D48").Value = "'124.28"
$ws.Range("E48").Value = "  +1.94%  "
$ws.Range("E49").Value = "  +1.91%  "
$ws.Range("E50").Value = "  +1.69%  "
$ws.Range("D51").Value = "'1.188"
$ws.Range("E51").Value = "  +3.24%  "

$ws.Range("D2:E51").Style = "Normal"
